# "Added a few more slots"
#
# The meta-description blurb (bold label "Meta description" + the
# descriptive sentence) that used to sit right under the H1 title is
# removed from there, and its descriptive sentence is moved down to the
# very end of the document, replacing the old "Create a fun cartoon
# image of a Maya warrior..." AI-art prompt paragraph. A new bold
# "Play Black and White Slot Game for Free" paragraph is inserted
# immediately above it.

$d = $word.ActiveDocument

# --- Step 1: turn the trailing "Maya warrior" image-prompt paragraph
# into the new closing pair of paragraphs ------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara  = $d.Paragraphs($lastIndex)

# Split off a brand-new empty paragraph right before it.
$lastPara.Range.InsertParagraphBefore()

$newParaIndex = $d.Paragraphs.Count - 1
$newPara      = $d.Paragraphs($newParaIndex)

# Populate that empty paragraph with a bold run via a raw OOXML
# fragment so the run formatting comes out exactly as authored
# (an empty leading run followed by the bold text run).
$boldFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Black and White Slot Game for Free</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($boldFragment)

# Now replace the text of the (still-last) Maya-warrior paragraph with
# the meta-description sentence, keeping its existing italic run.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalPara.Range.Find.Execute(
    "Create a fun cartoon image of a Maya warrior wearing glasses, with a big smile on their face. The warrior should be holding a slot machine handle in one hand, and surrounded by colorful flowers and butterflies. The background should feature a jungle landscape with a Mayan pyramid in the distance. Use a mix of bold colors to make the image pop and convey a sense of excitement and fun. The image should be eye-catching and playful, inviting players to explore the game and enjoy their gambling experience.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the medieval-themed Black and White slot game with unique black and white characters, free spins, and an expandable grid. Play for free now.",
    2)

# --- Step 2: drop the old "Meta description" paragraph from the top ---

$d.Paragraphs(2).Range.Delete()
